$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8 from 45175 to 45183.
$ws.Range("C2:C8").Value = 45183
